$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = 0
$ws.Range("H112").Value = 1625.9706
$ws.Range("J112").Value = 1706.1
$ws.Range("L112").Value = 5118.299999999999
$ws.Range("N112").Value = -7334.299999999999
$ws.Range("H116").Value = 721385.9
$ws.Range("I116").Value = 1113821.8
$ws.Range("K116").Value = 1113821.8
$ws.Range("M116").Value = -1110379.8
$ws.Range("H137").Value = 4363.93
$ws.Range("I137").Value = 4563.48
$ws.Range("K137").Value = 13690.44
$ws.Range("M137").Value = -11140.44
$ws.Range("H138").Value = 3929.6064
$ws.Range("I138").Value = 2156.1667
$ws.Range("J138").Value = 4189.1343
$ws.Range("K138").Value = 6468.500100000001
$ws.Range("L138").Value = 12567.4029
$ws.Range("M138").Value = -1328.500100000001
$ws.Range("N138").Value = -22847.4029

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 35197.75
$ws.Range("J23").Value = 41619.5
$ws.Range("L23").Value = 41619.5
$ws.Range("N23").Value = -42137.5
$ws.Range("H63").Value = 15394444
$ws.Range("J63").Value = 5211.4287
$ws.Range("L63").Value = 5211.4287
$ws.Range("N63").Value = -6583.4287
$ws.Range("H66").Value = 15394444
$ws.Range("J66").Value = 5211.4287
$ws.Range("L66").Value = 26057.1435
$ws.Range("N66").Value = -32921.14350000001
$ws.Range("H74").Value = 4756.3213
$ws.Range("I74").Value = 5032.7144
$ws.Range("J74").Value = 3927.1428
$ws.Range("K74").Value = 5032.7144
$ws.Range("L74").Value = 3927.1428
$ws.Range("M74").Value = -4158.7144
$ws.Range("N74").Value = -5675.1428
$ws.Range("H77").Value = 4756.3213
$ws.Range("I77").Value = 5032.7144
$ws.Range("J77").Value = 3927.1428
$ws.Range("K77").Value = 25163.572
$ws.Range("L77").Value = 19635.714
$ws.Range("M77").Value = -20795.572
$ws.Range("N77").Value = -28371.714
$ws.Range("H80").Value = 39979.145
$ws.Range("J80").Value = 39979.145
$ws.Range("L80").Value = 39979.145
$ws.Range("N80").Value = -41975.145
$ws.Range("H83").Value = 39979.145
$ws.Range("J83").Value = 39979.145
$ws.Range("L83").Value = 119937.435
$ws.Range("N83").Value = -129921.435
$ws.Range("H102").Value = 2127.9092
$ws.Range("I102").Value = 1901.3334
$ws.Range("J102").Value = 2399.8
$ws.Range("K102").Value = 1901.3334
$ws.Range("L102").Value = 2399.8
$ws.Range("M102").Value = -279.3334
$ws.Range("N102").Value = -5643.8
$ws.Range("H122").Value = 3036.8
$ws.Range("I122").Value = 1881
$ws.Range("J122").Value = 7660
$ws.Range("K122").Value = 5643
$ws.Range("L122").Value = 22980
$ws.Range("M122").Value = -3193
$ws.Range("N122").Value = -27880

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2679.2
$ws.Range("I86").Value = 2498.6667
$ws.Range("J86").Value = 2950
$ws.Range("K86").Value = 2498.6667
$ws.Range("L86").Value = 2950
$ws.Range("M86").Value = -1375.6667
$ws.Range("N86").Value = -5196
$ws.Range("H89").Value = 2679.2
$ws.Range("I89").Value = 2498.6667
$ws.Range("J89").Value = 2950
$ws.Range("K89").Value = 12493.3335
$ws.Range("L89").Value = 14750
$ws.Range("M89").Value = -6877.333500000001
$ws.Range("N89").Value = -25982
$ws.Range("H105").Value = 3400
$ws.Range("I105").Value = 2760
$ws.Range("K105").Value = 2760
$ws.Range("M105").Value = -1013
$ws.Range("H107").Value = 1653.875
$ws.Range("I107").Value = 1584.9
$ws.Range("J107").Value = 1768.8334
$ws.Range("K107").Value = 1584.9
$ws.Range("L107").Value = 1768.8334
$ws.Range("M107").Value = 335.0999999999999
$ws.Range("N107").Value = -5608.8334
$ws.Range("H134").Value = 3358.738
$ws.Range("I134").Value = 1775.9032
$ws.Range("J134").Value = 7819.4546
$ws.Range("K134").Value = 5327.7096
$ws.Range("L134").Value = 23458.3638
$ws.Range("M134").Value = -2792.7096
$ws.Range("N134").Value = -28528.3638

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4636.7617
$ws.Range("I31").Value = 2019.3846
$ws.Range("J31").Value = 5810.069
$ws.Range("K31").Value = 2019.3846
$ws.Range("L31").Value = 5810.069
$ws.Range("M31").Value = -1724.3846
$ws.Range("N31").Value = -6400.069
$ws.Range("H34").Value = 4636.7617
$ws.Range("I34").Value = 2019.3846
$ws.Range("J34").Value = 5810.069
$ws.Range("K34").Value = 2019.3846
$ws.Range("L34").Value = 5810.069
$ws.Range("M34").Value = -1817.3846
$ws.Range("N34").Value = -6214.069
$ws.Range("H68").Value = 46783.07
$ws.Range("J68").Value = 46783.07
$ws.Range("L68").Value = 46783.07
$ws.Range("N68").Value = -48281.07
$ws.Range("H71").Value = 46783.07
$ws.Range("J71").Value = 46783.07
$ws.Range("L71").Value = 140349.21
$ws.Range("N71").Value = -147837.21
$ws.Range("H87").Value = 19985.715
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 20753.846
$ws.Range("K87").Value = 10000
$ws.Range("L87").Value = 20753.846
$ws.Range("M87").Value = -8814
$ws.Range("N87").Value = -23125.846
$ws.Range("H90").Value = 19985.715
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 20753.846
$ws.Range("K90").Value = 30000
$ws.Range("L90").Value = 62261.538
$ws.Range("M90").Value = -24072
$ws.Range("N90").Value = -74117.538
$ws.Range("H123").Value = 38780
$ws.Range("J123").Value = 38780
$ws.Range("L123").Value = 38780
$ws.Range("N123").Value = -48580

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 212.17647
$ws.Range("J23").Value = 257.66666
$ws.Range("L23").Value = 772.9999799999999
$ws.Range("N23").Value = -1242.99998
$ws.Range("H33").Value = 198.70589
$ws.Range("I33").Value = 163.1
$ws.Range("J33").Value = 249.57143
$ws.Range("K33").Value = 978.5999999999999
$ws.Range("L33").Value = 1497.42858
$ws.Range("M33").Value = -695.5999999999999
$ws.Range("N33").Value = -2063.42858
$ws.Range("H80").Value = 3700.3809
$ws.Range("J80").Value = 3928.2222
$ws.Range("L80").Value = 11784.6666
$ws.Range("N80").Value = -13656.6666
$ws.Range("H83").Value = 3700.3809
$ws.Range("J83").Value = 3928.2222
$ws.Range("L83").Value = 35353.99980000001
$ws.Range("N83").Value = -44713.99980000001
$ws.Range("H97").Value = 508.22223
$ws.Range("J97").Value = 557.34784
$ws.Range("L97").Value = 1672.04352
$ws.Range("N97").Value = -2664.04352
$ws.Range("H113").Value = 579.2182
$ws.Range("I113").Value = 600.5806
$ws.Range("J113").Value = 551.625
$ws.Range("K113").Value = 1801.7418
$ws.Range("L113").Value = 1654.875
$ws.Range("M113").Value = 368.2582
$ws.Range("N113").Value = -5994.875
$ws.Range("H124").Value = 5600
$ws.Range("J124").Value = 5600
$ws.Range("L124").Value = 16800
$ws.Range("N124").Value = -26620

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716740
$ws.Range("I80").Value = 125001390
$ws.Range("K80").Value = 125001390
$ws.Range("M80").Value = -125000392
$ws.Range("H83").Value = 35716740
$ws.Range("I83").Value = 125001390
$ws.Range("K83").Value = 625006950
$ws.Range("M83").Value = -625001958
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H97").Value = 1519.4117
$ws.Range("I97").Value = 1055.3334
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 1055.3334
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -559.3334
$ws.Range("N97").Value = -5992

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4840.0835
$ws.Range("I132").Value = 2321.457
$ws.Range("J132").Value = 11621
$ws.Range("K132").Value = 6964.370999999999
$ws.Range("L132").Value = 34863
$ws.Range("M132").Value = -4434.370999999999
$ws.Range("N132").Value = -39923
